$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Update row 2 values (QR step matching Samsung devices)
$ws.Range("F2").Value = "QmFvMjQwNzE5OTc="
$ws.Range("A2").Value = "BAOTG2407"
$ws.Range("B2").Value = "'002704070016694"

# Change the selected/active cell (scroll to up page)
$ws.Range("G5").Select()
